# Apply the "nuevos experimentos no convexos" parameter update.
#
# The workbook stores every value (even the purely numeric-looking ones)
# as a shared string (t="s"), not as a number. Assigning a numeric-looking
# string straight to .Value would make Excel coerce it into a real number
# cell, which is not what we want here, so for values that parse as a
# number we go through .Formula with a leading apostrophe (forces text
# entry) and then reset the cell style back to "Normal" so no stray
# number-format / quote-prefix styling sticks around on the cell.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, [string]$addr, [string]$val)

    $looksNumeric = $val -match '^-?[0-9]+(\.[0-9]+)?([eE][-+]?[0-9]+)?$'

    if ($looksNumeric) {
        $ws.Range($addr).Formula = "'" + $val
        $ws.Range($addr).Style = "Normal"
    } else {
        $ws.Range($addr).Value = $val
    }
}

$wsLider = $wb.Worksheets.Item("Restricciones_del_lider")
Set-TextValue $wsLider "A2" "-0.9 + x"
Set-TextValue $wsLider "B2" "-0.09999999999999998"
Set-TextValue $wsLider "D2" "0.42"
Set-TextValue $wsLider "A3" "0.8999999999999999 - x"
Set-TextValue $wsLider "B3" "-1.9"
Set-TextValue $wsLider "D3" "0.02"

$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")
Set-TextValue $wsFollower "A2" "-2.7 + y"
Set-TextValue $wsFollower "B2" "1.7000000000000002"
Set-TextValue $wsFollower "D2" "0.29"
Set-TextValue $wsFollower "E2" "0"
Set-TextValue $wsFollower "F2" "0"
Set-TextValue $wsFollower "A3" "2.7 - y"
Set-TextValue $wsFollower "B3" "-3.7"
Set-TextValue $wsFollower "D3" "0.52"
Set-TextValue $wsFollower "E3" "0"
Set-TextValue $wsFollower "F3" "0"

$wsPunto = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $wsPunto "A2" "0.9"
Set-TextValue $wsPunto "B2" "2.7"

# NOTE: "Vector_bf" and "Vector_BF" are two distinct sheets that differ
# only by case; Worksheets.Item(<name>) resolves case-insensitively here
# and would hit the same ("Vector_bf") sheet for both, so address them by
# their 1-based position instead (5th and 6th tabs, respectively).
$wsBf = $wb.Worksheets.Item(5)
Set-TextValue $wsBf "A2" "-6.16"

$wsBF = $wb.Worksheets.Item(6)
Set-TextValue $wsBF "A2" "-1.7"
Set-TextValue $wsBF "A3" "-5.4"
